$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plans")

# Row 5: switch from Micronegocio / M_Con_TotalPlay_TV / 50 Megas
#        to     Residencial  / Con_TotalPlay_TV      / 20 Megas
$ws.Range("B5").Value = "Residencial"
$ws.Range("C5").Value = "Con_TotalPlay_TV"
$ws.Range("D5").Value = 20

# Row 6: switch from Micronegocio / M_Con_TotalPlay_TV / 120 Megas
#        to     Residencial  / Sin_TotalPlay_TV      / 50 Megas
$ws.Range("B6").Value = "Residencial"
$ws.Range("C6").Value = "Sin_TotalPlay_TV"
$ws.Range("D6").Value = 50

# Row 7: switch from Micronegocio / M_Con_TotalPlay_TV / 220 Megas
#        to     Residencial  / Sin_TotalPlay_TV      / 100 Megas
$ws.Range("B7").Value = "Residencial"
$ws.Range("C7").Value = "Sin_TotalPlay_TV"
$ws.Range("D7").Value = 100

# Row 8: switch from Micronegocio / M_Con_TotalPlay_TV / 520 Megas
#        to     Residencial  / Sin_TotalPlay_TV      / 500 Megas
$ws.Range("B8").Value = "Residencial"
$ws.Range("C8").Value = "Sin_TotalPlay_TV"
$ws.Range("D8").Value = 500

# Row 9: switch from Micronegocio / M_Con_TotalPlay_TV / 1000 Megas
#        to     Residencial  / Sin_TotalPlay_TV      / 1000 Megas (unchanged)
$ws.Range("B9").Value = "Residencial"
$ws.Range("C9").Value = "Sin_TotalPlay_TV"
$ws.Range("D9").Value = 1000

$ws.Range("E15").Select()
